# Update cryptos list (price + volume change) per the source commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.620.31'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.01%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.61'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.55%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.28'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.41%  '

# Row 6
$ws.Range("E6").Value = '  +0.45%  '

# Row 7
$ws.Range("E7").Value = '  +0.40%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3909'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.21%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07883'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9680'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.27'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.825.26'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.715'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.40%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.917'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.31%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06939'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.06%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.39'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.86%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001004'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.92'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.06%  '

# Row 20
$ws.Range("E20").Value = '  +0.21%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.619.31'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.97%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.313'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.63%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.04'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.121'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.50%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.056.78'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.07%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.25'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.744'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.41%  '

# Row 29
$ws.Range("E29").Value = '  +0.32%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.02'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.84%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09345'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9331'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.42%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.310'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.19%  '

# Row 34
$ws.Range("E34").Value = '  +0.55%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.346'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.85%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05826'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.55%  '

# Row 37
$ws.Range("E37").Value = '  -1.87%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.154'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.14%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.898'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.40%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5644'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.17%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.914'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.19%  '

# Row 42
$ws.Range("E42").Value = '  -0.52%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07250'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.98%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.191'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.43%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.66'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.34%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5308'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("E47").Value = '  -8.24%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.847'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.02%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.33'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.02%  '

# Row 50
$ws.Range("E50").Value = '  +0.46%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.344'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.85%  '
